# edit.ps1 - apply the "What's That Requirements Document" edits:
#  1) In the "image / audio capture tab" paragraph:
#       a. split the leading "        <TAB>" run so the TAB itself is no
#          longer bold (only the 8 leading spaces stay bold).
#       b. rewrite the descriptive sentence about the icon behaviour.
#  2) In the "boards tab" paragraph, extend the sentence about filters and
#     bold just the new trailing period.

$d = $word.ActiveDocument

$RSQUO = [char]0x2019

# ---------------------------------------------------------------------
# 1a) Un-bold the tab character that precedes
#     "The image / audio capture tab allows users..."
# ---------------------------------------------------------------------
$paras = $d.Paragraphs
$count = $paras.Count
for ($i = 1; $i -le $count; $i++) {
    $cand = $paras.Item($i)
    if ($cand.Range.Text -like "*The image / audio capture tab allows users*") {
        $targetPara = $cand
        break
    }
}

$pStart = $targetPara.Range.Start
$tabRng = $d.Range($pStart + 8, $pStart + 9)
if ($tabRng.Text -eq [char]9) {
    $tabRng.Bold = 0
}

# ---------------------------------------------------------------------
# 1b) Replace the icon/location sentence with the expanded description.
# ---------------------------------------------------------------------
$oldIcon = "An icon with the object" + $RSQUO + "s name and location will be prompted on the screen; if the user clicks on the icon, more information will be presented to the user."
$newIcon = "An icon will pop up displaying the specific name of the object (Ex - a M1 MacBook Air). If the image associated happens to have a location (Example - Empire State Building), information on the image" + $RSQUO + "s location will also be displayed. If the user clicks on the icon, more information will be presented to the user."

$found = $d.Content.Find.Execute($oldIcon, $false, $false, $false, $false, $false, $true, 1, $false, $newIcon, 2)
if (-not $found) {
    throw "Could not find the icon/location sentence to replace."
}

# ---------------------------------------------------------------------
# 2) Extend the boards-tab sentence and bold the new trailing period.
# ---------------------------------------------------------------------
$oldBoards = "The boards tab will include all of the community boards. Users can apply filters to see which boards are trending near their area, boards that are trending worldwide, etc. Users can also search for boards they desire on the search bar. On clicking a board group, the user will be shown all photos within that board. Users can then be able to, like, leave a comment, or flag a post."
$newBoards = "The boards tab will include all of the community boards. Users can apply filters to see which boards are trending near their area, boards that are trending worldwide, and boards that are gaining high attraction. Users can also search for boards they desire on the search bar. On clicking a board group, the user will be shown all photos within that board. Users can then be able to, like, leave a comment, or flag a post."

$found2 = $d.Content.Find.Execute($oldBoards, $false, $false, $false, $false, $false, $true, 1, $false, $newBoards, 2)
if (-not $found2) {
    throw "Could not find the boards-tab sentence to replace."
}

# Bold just the period right after "...gaining high attraction"
$rng = $d.Content.Duplicate
$found3 = $rng.Find.Execute("attraction.", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found3) {
    throw "Could not find the 'attraction.' anchor to bold the period."
}
$periodRng = $d.Range($rng.End - 1, $rng.End)
if ($periodRng.Text -eq ".") {
    $periodRng.Bold = 1
}
